$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Switch focus to the ProductLoanInput sheet (it becomes the active/selected tab)
$wsInput.Activate()

# Update the repayment strategy value in B17 (row 17, "repaymentstrategy")
$wsInput.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Leave the selection on the cell that was edited
$wsInput.Range("B17").Select()
